# Add data for 2021-12-18: update "through 12-09" snapshot to "through 12-10"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and workbook title to reflect the new "through" date
$ws.Name = "Through 2021-12-10"

# Row 12 (October) - 2021 column group (T/U/V)
$ws.Range("U12").Value = 188
$ws.Range("V12").Value = 0.0408

# Row 14 (December) - label update
$ws.Range("A14").Value = "December (through 12-10)"

# 2015 group (B/C/D)
$ws.Range("C14").Value = 7
$ws.Range("D14").Value = 0.2222

# 2016 group (E/F/G)
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 25
$ws.Range("G14").Value = 0.1071

# 2018 group (K/L/M)
$ws.Range("L14").Value = 22
$ws.Range("M14").Value = 0.0833

# 2019 group (N/O/P)
$ws.Range("O14").Value = 11
$ws.Range("P14").Value = 0.2143

# 2020 group (Q/R/S)
$ws.Range("R14").Value = 46
$ws.Range("S14").Value = 0.0417

# 2021 group (T/U) - total only column, no arrest_rate for this row
$ws.Range("U14").Value = 75

# Row 15 (Total)
$ws.Range("C15").Value = 265
$ws.Range("D15").Value = 0.1167

$ws.Range("E15").Value = 63
$ws.Range("F15").Value = 528
$ws.Range("G15").Value = 0.1066

$ws.Range("L15").Value = 630
$ws.Range("M15").Value = 0.1076

$ws.Range("O15").Value = 491
$ws.Range("P15").Value = 0.104

$ws.Range("R15").Value = 1246
$ws.Range("S15").Value = 0.0503

$ws.Range("U15").Value = 1620
$ws.Range("V15").Value = 0.0576
